$wb = $excel.ActiveWorkbook

$hbUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bd07102de3937c6688ac1e30a2230b5f96fb0ce/e2e/7d818bb1-04de-4d02-88e9-033dccb47dd1.md"
$hbName = "7d818bb1-04de-4d02-88e9-033dccb47dd1.md"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("I7").Value = $hbName
$ws.Range("J7").Value = "7d818bb1-04de-4d02-88e9-033dccb47dd1.7d4ed1c3ab68396f358067124cc98e75e2223fd0.zh-cn.xlf"
$ws.Range("K7").Value = "2016-08-28 20:54:46"
$ws.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49d794a0240a5d4ab5da565918cbec3907672014/e2e/7d818bb1-04de-4d02-88e9-033dccb47dd1.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bd07102de3937c6688ac1e30a2230b5f96fb0ce/e2e/7d818bb1-04de-4d02-88e9-033dccb47dd1.md."
$ws.Hyperlinks.Add($ws.Range("I7"), $hbUrl, [System.Type]::Missing, [System.Type]::Missing, $hbName)

# ---- de-de sheet ----
$ws2 = $wb.Worksheets.Item("de-de")
$ws2.Range("I7").Value = $hbName
$ws2.Range("J7").Value = "7d818bb1-04de-4d02-88e9-033dccb47dd1.7d4ed1c3ab68396f358067124cc98e75e2223fd0.de-de.xlf"
$ws2.Range("K7").Value = "2016-08-28 20:54:52"
$ws2.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49d794a0240a5d4ab5da565918cbec3907672014/e2e/7d818bb1-04de-4d02-88e9-033dccb47dd1.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bd07102de3937c6688ac1e30a2230b5f96fb0ce/e2e/7d818bb1-04de-4d02-88e9-033dccb47dd1.md."
$ws2.Hyperlinks.Add($ws2.Range("I7"), $hbUrl, [System.Type]::Missing, [System.Type]::Missing, $hbName)
